# Rule_016 "User added to Microsoft Entra ID Privileged Groups" triage template
# Populate the "Output" column (F) for rows 3-12 with sample KQL output tables,
# resize rows to fit the pasted multi-line content, tighten column C's width, and
# leave the view scrolled/selected where the analyst left off (F3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text37 = @"
TimeGenerated            AssignedUser              RoleName                   Initiator              Result
2025-10-09 10:45:58     emma.wilson@abc.com       Global Administrator       it.admin@abc.com       Success
2025-10-09 08:45:58     david.brown@abc.com       Privileged Role Admin      security.lead@abc.com  Success  
2025-10-09 06:45:58     lisa.garcia@abc.com       Security Administrator     it.admin@abc.com       Success
2025-10-09 03:45:58     robert.jones@abc.com      Exchange Administrator     exchange.admin@abc.com Success
"@

$text38 = @"
TimeGenerated            Username                  GroupName                  AddedBy                IPAddress     Location
2025-10-09 10:45:58     emma.wilson@abc.com       Global Administrators      it.admin@abc.com       172.16.1.100  New York, NY, US
2025-10-09 08:45:58     david.brown@abc.com       Privileged Role Admins     security.lead@abc.com  172.16.2.50   Chicago, IL, US
2025-10-09 06:45:58     lisa.garcia@abc.com       Security Administrators    it.admin@abc.com       172.16.1.100  New York, NY, US
2025-10-09 03:45:58     robert.jones@abc.com      Exchange Administrators    exchange.admin@abc.com 172.16.3.75   Dallas, TX, US
"@

$text39 = @"
UserPrincipalName        AssignedRole               RiskLevel  RequiresApproval  BusinessJustification                    IsHighRisk
emma.wilson@abc.com      Global Administrator       Critical   True              Temporary admin access for migration    True
david.brown@abc.com      Privileged Role Admin      High       True              PIM role for security operations        True
lisa.garcia@abc.com      Security Administrator     Medium     False             Standard security analyst role          False
robert.jones@abc.com     Exchange Administrator     Medium     False             Exchange maintenance duties             False
"@

$text40 = @"
UserPrincipalName        SignInCount  UniqueIPs  UniqueLocations  FailedSignIns  Locations                    IPs                              RiskScore
emma.wilson@abc.com      45          3          2                2              New York, NY; Boston, MA    172.16.1.100; 203.0.113.25      12
david.brown@abc.com      32          2          1                0              Chicago, IL                  172.16.2.50; 10.0.1.25          6
lisa.garcia@abc.com      67          1          1                1              New York, NY                 172.16.1.100                     8
robert.jones@abc.com     28          4          3                5              Dallas, TX; Austin, TX       172.16.3.75; 192.0.2.100        27
"@

$text41 = @"
TimeGenerated            InitiatedBy              OperationName             TargetUser              InitiatorRole            SourceIP      AuthenticationMethod     Result
2025-10-09 10:45:58     it.admin@abc.com         Add member to role        emma.wilson@abc.com     Global Administrator     172.16.1.100  MFA + Smart Card        Success
2025-10-09 08:45:58     security.lead@abc.com    Add eligible member       david.brown@abc.com     Privileged Role Admin    172.16.2.50   MFA + Authenticator     Success
2025-10-09 06:45:58     it.admin@abc.com         Add member to role        lisa.garcia@abc.com     Global Administrator     172.16.1.100  MFA + Smart Card        Success
2025-10-09 03:45:58     exchange.admin@abc.com   Add member to role        robert.jones@abc.com    Exchange Administrator   172.16.3.75   Password Only           Success
"@

$text42 = @"
UserPrincipalName        SuspiciousActivity                    ThreatIndicators              VirusTotalResult                             RequiresScreenshot  ThreatLevel
emma.wilson@abc.com      Multiple location sign-ins in 1hr    Impossible travel detected    Clean - No malicious indicators              True               Medium
david.brown@abc.com      No suspicious activity detected       None                          Clean - No malicious indicators              False              Low
robert.jones@abc.com     Failed sign-ins from unknown IPs     Brute force attempt           Flagged - 2 vendors marked IP suspicious    True               High
"@

$text43 = @"
UserPrincipalName        ActionTaken                          SessionsRevoked  MFAStatus                  ITNotified  IAMNotified  ComplianceStatus
emma.wilson@abc.com      Account review initiated             0               Enabled - Authenticator    True        True         Under Review
david.brown@abc.com      No action required                   0               Enabled - SMS + App        False       False        Compliant
lisa.garcia@abc.com      Standard monitoring                  0               Enabled - Smart Card       False       False        Compliant
robert.jones@abc.com     Session revoked, password reset      5               Disabled - Requires setup  True        True         Non-Compliant
"@

$text44 = @"
InvestigationID  UserPrincipalName        ITContactStatus                    IAMVerificationStatus                CredentialResetRequired  FollowUpAction
PRIV-2025-001    emma.wilson@abc.com      Contacted - Awaiting verification  Verified - Legitimate business need  False                   Monitor for 48 hours
PRIV-2025-002    david.brown@abc.com      Not required                       Auto-approved - Standard role        False                   No follow-up required
PRIV-2025-003    lisa.garcia@abc.com      Not required                       Auto-approved - Standard role        False                   No follow-up required
PRIV-2025-004    robert.jones@abc.com     Contacted - Credential reset       Failed verification - Unauthorized   True                    Full account audit initiated
"@

$text45 = @"
InvestigationID  UserAffected             FindingsSummary                                      RemediationActions                                     Status                   AssignedAnalyst
PRIV-2025-001    emma.wilson@abc.com      Legitimate assignment requiring monitoring           Enable additional monitoring, require approval        In Progress - Monitoring L2-Analyst-02
PRIV-2025-002    david.brown@abc.com      Standard privileged role - no issues found          No action required - maintain current access          Closed - No Action       L1-Analyst-04
PRIV-2025-003    lisa.garcia@abc.com      Security role assignment approved and compliant     No action required - maintain current access          Closed - No Action       L1-Analyst-05
PRIV-2025-004    robert.jones@abc.com     Unauthorized privileged access - critical finding   Remove access, reset credentials, enable strict CAP   Active - Critical        L3-Analyst-01
"@

# Column C was manually narrowed after the Output column was filled in.

$ws.Columns.Item(3).ColumnWidth = 32.95

$ws.Range("F3").Value = $text37
$ws.Rows.Item(3).RowHeight = 345.6

$ws.Range("F4").Value = $text37
$ws.Rows.Item(4).RowHeight = 345.6

$ws.Range("F5").Value = $text38
$ws.Rows.Item(5).RowHeight = 388.8

$ws.Range("F6").Value = $text39
$ws.Rows.Item(6).RowHeight = 345.6

$ws.Range("F7").Value = $text40
$ws.Rows.Item(7).RowHeight = 345.6

$ws.Range("F8").Value = $text41
$ws.Rows.Item(8).RowHeight = 409.6

$ws.Range("F9").Value = $text42
$ws.Rows.Item(9).RowHeight = 345.6

$ws.Range("F10").Value = $text43
$ws.Rows.Item(10).RowHeight = 360

$ws.Range("F11").Value = $text44
$ws.Rows.Item(11).RowHeight = 409.6

$ws.Range("F12").Value = $text45
$ws.Rows.Item(12).RowHeight = 409.6


# Final view state saved with the workbook: frozen pane scrolled to the top,
# active cell on the first populated Output cell.
$ws.Range("F3").Select() | Out-Null
